$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts existing rows 6-39 down to 7-40)
$ws.Rows("6:6").Insert()

# Populate the new row with the "play" / "PLAY" key-value pair
$ws.Range("A6").Value = "play"
$ws.Range("B6").Value = "PLAY"

# Update the selection to match the new active cell (A6)
$null = $ws.Range("A6").Select()
